$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 649 entirely (the "bear" post); all rows below shift up by one.
$ws.Rows(649).Delete()
